# Actualización automática 2025-10-20 09:30:08
# Inserts a new client row ("CARANGUI ARMIJOS LUIS FRANCISCO") at row 14
# on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing the
# existing alphabetically-sorted client rows down by one. The new row gets
# zero values for every metric column (the client has no sales yet).
# The trailing summary rows are also refreshed: the "X de 53" counters on
# "VENTAS POR GRUPO" become "X de 54" to reflect the new client count.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": columns A:R, data rows 2-54 (was 2-54 -> 55),
# header row 1, trailing "X de N" summary row.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row before the old row 14 (CEDEÑO MACIAS...), shifting
# every row from 14 down to 55 (and the dimension) down by one.
$ws1.Rows.Item(14).Insert()

$ws1.Cells.Item(14, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws1.Cells.Item(14, 2).Value = "CARANGUI ARMIJOS LUIS FRANCISCO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(14, $c).Value = 0
}

# The final summary row (previously row 55, now row 56) reports "X de 53";
# bump it to "X de 54" now that there are 54 clients listed.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(56, $c)
    $text = $cell.Value()
    $cell.Value = ($text -replace "53", "54")
}

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": columns A:G, data rows 2-58 (was 2-58 -> 59),
# header row 1, trailing numeric totals row (values unaffected, only its
# row number shifts down with everything else).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(14).Insert()

$ws2.Cells.Item(14, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws2.Cells.Item(14, 2).Value = "CARANGUI ARMIJOS LUIS FRANCISCO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(14, $c).Value = 0
}
